$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-03 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-04 Thursday", 2) | Out-Null
$d.Content.Find.Execute("977÷9=108, 5", $true, $false, $false, $false, $false, $true, 1, $false, "527÷4=131, 3", 2) | Out-Null
$d.Content.Find.Execute("724÷6=120, 4", $true, $false, $false, $false, $false, $true, 1, $false, "432÷6=72, 0", 2) | Out-Null
$d.Content.Find.Execute("219÷2=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "309÷8=38, 5", 2) | Out-Null
$d.Content.Find.Execute("557÷3=185, 2", $true, $false, $false, $false, $false, $true, 1, $false, "852÷6=142, 0", 2) | Out-Null
$d.Content.Find.Execute("645÷4=161, 1", $true, $false, $false, $false, $false, $true, 1, $false, "290÷4=72, 2", 2) | Out-Null
$d.Content.Find.Execute("351÷6=58, 3", $true, $false, $false, $false, $false, $true, 1, $false, "406÷7=58, 0", 2) | Out-Null
$d.Content.Find.Execute("817÷6=136, 1", $true, $false, $false, $false, $false, $true, 1, $false, "799÷2=399, 1", 2) | Out-Null
$d.Content.Find.Execute("120÷8=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "403÷5=80, 3", 2) | Out-Null
$d.Content.Find.Execute("721÷7=103, 0", $true, $false, $false, $false, $false, $true, 1, $false, "963÷7=137, 4", 2) | Out-Null
$d.Content.Find.Execute("935÷3=311, 2", $true, $false, $false, $false, $false, $true, 1, $false, "277÷6=46, 1", 2) | Out-Null
$d.Content.Find.Execute("293÷7=41, 6", $true, $false, $false, $false, $false, $true, 1, $false, "430÷6=71, 4", 2) | Out-Null
$d.Content.Find.Execute("493÷3=164, 1", $true, $false, $false, $false, $false, $true, 1, $false, "131÷6=21, 5", 2) | Out-Null
$d.Content.Find.Execute("196÷2=98, 0", $true, $false, $false, $false, $false, $true, 1, $false, "881÷6=146, 5", 2) | Out-Null
$d.Content.Find.Execute("860÷8=107, 4", $true, $false, $false, $false, $false, $true, 1, $false, "347÷6=57, 5", 2) | Out-Null
$d.Content.Find.Execute("911÷4=227, 3", $true, $false, $false, $false, $false, $true, 1, $false, "247÷7=35, 2", 2) | Out-Null
$d.Content.Find.Execute("815÷5=163, 0", $true, $false, $false, $false, $false, $true, 1, $false, "634÷8=79, 2", 2) | Out-Null
$d.Content.Find.Execute("536÷3=178, 2", $true, $false, $false, $false, $false, $true, 1, $false, "418÷4=104, 2", 2) | Out-Null
$d.Content.Find.Execute("857÷5=171, 2", $true, $false, $false, $false, $false, $true, 1, $false, "112÷9=12, 4", 2) | Out-Null
$d.Content.Find.Execute("796÷4=199, 0", $true, $false, $false, $false, $false, $true, 1, $false, "585÷9=65, 0", 2) | Out-Null
$d.Content.Find.Execute("816÷2=408, 0", $true, $false, $false, $false, $false, $true, 1, $false, "524÷5=104, 4", 2) | Out-Null
$d.Content.Find.Execute("932÷4=233, 0", $true, $false, $false, $false, $false, $true, 1, $false, "899÷9=99, 8", 2) | Out-Null
$d.Content.Find.Execute("333÷2=166, 1", $true, $false, $false, $false, $false, $true, 1, $false, "568÷3=189, 1", 2) | Out-Null
$d.Content.Find.Execute("666÷3=222, 0", $true, $false, $false, $false, $false, $true, 1, $false, "937÷7=133, 6", 2) | Out-Null
$d.Content.Find.Execute("730÷7=104, 2", $true, $false, $false, $false, $false, $true, 1, $false, "523÷5=104, 3", 2) | Out-Null
$d.Content.Find.Execute("959÷6=159, 5", $true, $false, $false, $false, $false, $true, 1, $false, "275÷3=91, 2", 2) | Out-Null
